# Fix Training Data Issue: the BF column held the source-file's date stamp
# ("5-7-2011-12", from the original "5-7-2011-12.xlsx" filename) instead of
# the actual game date. NBA.com stats display the "prior night" date for
# games, so the true calendar date is one day later: 2012-05-07.
#
# We write the corrected date as a literal text value (not an Excel date
# serial) to match the original inlineStr/text cell type. Assigning a
# formula like ="2012-05-07" and then Copy + PasteSpecial(xlPasteValues)
# converts it to a plain string value without Excel's automatic
# text->date parsing and without introducing any new cell style/format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("BF2:BF31")
$range.Formula = '="2012-05-07"'
$range.Copy()
$range.PasteSpecial(-4163)   # xlPasteValues
